$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Columns("F").Insert()
$ws.Rows("1").AutoFit()
$ws.Range("F1").Value = "LOCATION"
$ws.Columns("F").ColumnWidth = 10.33
$null = $ws.Range("C1").Select()
